$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores prices/percent-changes as plain text (e.g. "60.942.41",
# "  -0.16%  "). Force columns D:E to Text format before writing so that
# numeric-looking strings (e.g. "592.78") are kept verbatim as text instead of
# being auto-converted to numbers by Excel's input parser.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 40 / Row 41: Kaspa and Stacks swapped position in the ranking.
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "2.03"
$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").Value = "  -0.42%  "

# Updated prices / 1h volume change percentages.
$ws.Range("D2").Value = "61.383.81"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.937.25"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "592.78"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "146.71"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "2.930.81"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "6.82"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "0.146"
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "33.92"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D16").Value = "3.404.40"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "61.204.12"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "6.78"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "2.907.27"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").Value = "433.95"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "13.61"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "0.687"
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").Value = "7.13"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").Value = "81.19"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "11.05"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "2.25"
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("D27").Value = "12.18"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "2.36"
$ws.Range("E29").Value = "  +7.91%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "2.63"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "7.18"
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").Value = "26.79"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "0.109"
$ws.Range("E34").Value = "  +1.86%  "
$sub3 = [string][char]0x2083
$ws.Range("D35").Value = "0.0" + $sub3 + "0873"
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").Value = "1.02"
$ws.Range("D37").Value = "3.14"
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("D38").Value = "5.67"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "49.91"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D42").Value = "8.67"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").Value = "0.295"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "40.50"
$ws.Range("E44").Value = "  -4.80%  "
$ws.Range("D45").Value = "383.63"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").Value = "0.0351"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "2.715.01"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").Value = "129.92"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("D50").Value = "24.38"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("E51").Value = "  +0.71%  "
